$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 9,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 5.682516
$arr[0,3] = 17.047548
$arr[0,4] = 0.4522589164991918
$arr[0,5] = 0.4522589164991919
$arr[0,6] = 1
$arr[0,7] = 0.3333333333333333
$arr[0,8] = 0.002213333333333333
$arr[0,9] = 0.00664
$arr[0,10] = 0.03176273618751495
$arr[0,11] = 0.03176273618751495
$arr[0,12] = 0.01257730208
$arr[0,13] = 0.11319571872
$arr[0,14] = 0.01436498065321518
$arr[0,15] = 0.01436498065321518
$arr[1,0] = 3
$arr[1,1] = 1
$arr[1,2] = 5.682516
$arr[1,3] = 17.047548
$arr[1,4] = 0.4522589164991918
$arr[1,5] = 0.4522589164991919
$arr[1,6] = 2
$arr[1,7] = 0.6666666666666666
$arr[1,8] = 0.06286566666666667
$arr[1,9] = 0.188597
$arr[1,10] = 0.9021621621621622
$arr[1,11] = 0.9021621621621622
$arr[1,12] = 0.357235156684
$arr[1,13] = 3.215116410156
$arr[1,14] = 0.4080108819660276
$arr[1,15] = 0.4080108819660277
$arr[2,0] = 3
$arr[2,1] = 1
$arr[2,2] = 5.682516
$arr[2,3] = 17.047548
$arr[2,4] = 0.4522589164991918
$arr[2,5] = 0.4522589164991919
$arr[2,6] = 1
$arr[2,7] = 0.3333333333333333
$arr[2,8] = 0.004604333333333334
$arr[2,9] = 0.013813
$arr[2,10] = 0.06607510165032289
$arr[2,11] = 0.06607510165032289
$arr[2,12] = 0.026164197836
$arr[2,13] = 0.235477780524
$arr[2,14] = 0.02988305387994899
$arr[2,15] = 0.029883053879949
$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 5.449245
$arr[3,3] = 16.347735
$arr[3,4] = 0.4336933920535619
$arr[3,5] = 0.433693392053562
$arr[3,6] = 1
$arr[3,7] = 0.3333333333333333
$arr[3,8] = 0.002213333333333333
$arr[3,9] = 0.00664
$arr[3,10] = 0.03176273618751495
$arr[3,11] = 0.03176273618751495
$arr[3,12] = 0.0120609956
$arr[3,13] = 0.1085489604
$arr[3,14] = 0.01377528879806578
$arr[3,15] = 0.01377528879806578
$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 5.449245
$arr[4,3] = 16.347735
$arr[4,4] = 0.4336933920535619
$arr[4,5] = 0.433693392053562
$arr[4,6] = 2
$arr[4,7] = 0.6666666666666666
$arr[4,8] = 0.06286566666666667
$arr[4,9] = 0.188597
$arr[4,10] = 0.9021621621621622
$arr[4,11] = 0.9021621621621622
$arr[4,12] = 0.342570419755
$arr[4,13] = 3.083133777795
$arr[4,14] = 0.3912617682904837
$arr[4,15] = 0.3912617682904838
$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 5.449245
$arr[5,3] = 16.347735
$arr[5,4] = 0.4336933920535619
$arr[5,5] = 0.433693392053562
$arr[5,6] = 1
$arr[5,7] = 0.3333333333333333
$arr[5,8] = 0.004604333333333334
$arr[5,9] = 0.013813
$arr[5,10] = 0.06607510165032289
$arr[5,11] = 0.06607510165032289
$arr[5,12] = 0.025090140395
$arr[5,13] = 0.225811263555
$arr[5,14] = 0.02865633496501244
$arr[5,15] = 0.02865633496501245
$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 1.432979666666667
$arr[6,3] = 4.298939
$arr[6,4] = 0.1140476914472462
$arr[6,5] = 0.1140476914472462
$arr[6,6] = 1
$arr[6,7] = 0.3333333333333333
$arr[6,8] = 0.002213333333333333
$arr[6,9] = 0.00664
$arr[6,10] = 0.03176273618751495
$arr[6,11] = 0.03176273618751495
$arr[6,12] = 0.003171661662222222
$arr[6,13] = 0.02854495496
$arr[6,14] = 0.003622466736233986
$arr[6,15] = 0.003622466736233987
$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 1.432979666666667
$arr[7,3] = 4.298939
$arr[7,4] = 0.1140476914472462
$arr[7,5] = 0.1140476914472462
$arr[7,6] = 2
$arr[7,7] = 0.6666666666666666
$arr[7,8] = 0.06286566666666667
$arr[7,9] = 0.188597
$arr[7,10] = 0.9021621621621622
$arr[7,11] = 0.9021621621621622
$arr[7,12] = 0.09008522206477777
$arr[7,13] = 0.8107669985830001
$arr[7,14] = 0.1028895119056508
$arr[7,15] = 0.1028895119056508
$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 1.432979666666667
$arr[8,3] = 4.298939
$arr[8,4] = 0.1140476914472462
$arr[8,5] = 0.1140476914472462
$arr[8,6] = 1
$arr[8,7] = 0.3333333333333333
$arr[8,8] = 0.004604333333333334
$arr[8,9] = 0.013813
$arr[8,10] = 0.06607510165032289
$arr[8,11] = 0.06607510165032289
$arr[8,12] = 0.006597916045222222
$arr[8,13] = 0.059381244407
$arr[8,14] = 0.007535712805361453
$arr[8,15] = 0.007535712805361455

$ws.Range("E2:T10").Value = $arr
